$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text updates (C3 / D3 headers for the std-deviation row) ---
$ws.Range("C3").Value = "DesvioPadrao F"
$ws.Range("D3").Value = "Desvio Padrão l"

# --- Bold + theme colour on the "Media" (average) and "Desvio Padrão" (stdev) result cells ---
$ws.Range("C2:D2").Font.Bold = $true
$ws.Range("C2:D2").Font.ThemeColor = 3
$ws.Range("C4:D4").Font.Bold = $true
$ws.Range("C4:D4").Font.ThemeColor = 3

# --- D1 header gains a left border (full box outline) ---
$ws.Range("D1").Borders.LineStyle = 1

# --- Remove the now-unused helper cells in rows 9-12 (columns A & B) ---
$ws.Range("A9:B12").Clear()

# --- Column widths: widen A,B and add width to (empty) E,F to match the new layout ---
$ws.Columns("A:B").ColumnWidth = 13.86
$ws.Columns("E:F").ColumnWidth = 13.86

# --- Move the scatter chart down/right (same size, new anchor) ---
$co = $ws.ChartObjects(1)
$co.Left = 15.75
$co.Top = 197.25

# --- Update the active selection ---
$ws.Range("D10").Select()
